$d = $word.ActiveDocument

function New-WordXmlPackage([string]$innerBody) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $innerBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# Step 1: the "Program executes ... no memory leaks" paragraph (originally
# paragraph 7) gets its split "leaks "/":"/" Yes" runs collapsed into two
# plain runs, and five further checklist paragraphs (release executable,
# project files, archive, estimate hours, submission count) are inserted
# right after it in simplified form, with the "submission count" answer
# changed from "Once" to "Twice".
# ---------------------------------------------------------------------------
$inner1 = '<w:p><w:r><w:t xml:space="preserve">Program executes without crashing Program has no memory leaks </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>: Yes</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t xml:space="preserve">A release executable has been made and included in the submission </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>: Yes</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t xml:space="preserve">Project files and source code are included in the submission </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>: Yes</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t xml:space="preserve">All files are packaged in a single compressed archive </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>: Yes</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t xml:space="preserve">Estimate the number of hours taken to complete this assessment </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>: 16</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>How many times have you submitted this assessment (including this time)?</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Twice</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'

$memoryLeaksPara = $d.Paragraphs.Item(7)
[void]$memoryLeaksPara.Range.InsertXML((New-WordXmlPackage $inner1))

# ---------------------------------------------------------------------------
# Step 2: the pathfinding-algorithm paragraph (the original paragraph 13,
# now paragraph 18 after the five-paragraph insert above) gets its bold
# "File:"/"NodeGraph.h"/"Line number:"/"61" runs merged into a single bold
# run, dropping the spell-check proofErr markers.
# ---------------------------------------------------------------------------
$innerPathfinding = '<w:p><w:r><w:t xml:space="preserve">The program implements a pathfinding </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">algorithm </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Yes, File: NodeGraph.h, Line number: 61</w:t></w:r></w:p>'
$pathfindingPara = $d.Paragraphs.Item(18)
[void]$pathfindingPara.Range.InsertXML((New-WordXmlPackage $innerPathfinding))

# ---------------------------------------------------------------------------
# Step 3: the NPC-AI-strategy paragraph (the original paragraph 14, now
# paragraph 19) gets the same "File:"/"Guard.h"/"Line number:"/"27" run
# merge and proofErr removal.
# ---------------------------------------------------------------------------
$innerNpcAi = '<w:p><w:r><w:t xml:space="preserve">The program implements an NPC AI </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">strategy </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Yes, File: Guard.h, Line number: 27</w:t></w:r></w:p>'
$npcAiPara = $d.Paragraphs.Item(19)
[void]$npcAiPara.Range.InsertXML((New-WordXmlPackage $innerNpcAi))

# ---------------------------------------------------------------------------
# Step 4: the original (now duplicated) "release executable" / "project
# files" / "archive" / "estimate hours" / "submission count (Once)"
# paragraphs - now sitting at positions 13-17, right after the new
# simplified copies and right before the pathfinding paragraph - are
# removed outright.
# ---------------------------------------------------------------------------
$firstDupe = $d.Paragraphs.Item(13)
$lastDupe = $d.Paragraphs.Item(17)
$dupeRange = $d.Range($firstDupe.Range.Start, $lastDupe.Range.End)
$dupeRange.Delete()
